$d = $word.ActiveDocument

# --- Body -------------------------------------------------------------
# The whole "What is Functional Programming?" write-up (intro paragraph,
# the odd-number example with its picture, and the Imperative/Declarative
# discussion) is removed, leaving a single trailing empty paragraph.

# 1) Drop the first paragraph ("What is Functional Programming ?")
#    entirely, mark and all.
$d.Paragraphs.Item(1).Range.Delete()

# 2) The paragraph that is now first ("It is a way of programming...")
#    keeps its paragraph mark/formatting, but loses its text.
$p1 = $d.Paragraphs.Item(1)
$keepRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$keepRange.Delete()

# 3) Remove every paragraph between this now-empty paragraph and the
#    document's final paragraph mark (the odd-number example, the
#    embedded picture, and the rest of the discussion).
$lastIndex = $d.Paragraphs.Count
if ($lastIndex -gt 2) {
    $start = $d.Paragraphs.Item(2).Range.Start
    $end = $d.Paragraphs.Item($lastIndex - 1).Range.End
    $d.Range($start, $end).Delete()
}

# 4) Merge the now-empty leading paragraph away so the document body
#    ends with exactly one empty paragraph before the section break.
if ($d.Paragraphs.Count -gt 1) {
    $d.Paragraphs.Item(1).Range.Delete()
}

# --- Header -------------------------------------------------------------
# "Functional programming" -> "optional"
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdr.Range.Find.Execute("Functional programming", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "optional", 2)
